$d = $word.ActiveDocument

# The document has these paragraphs:
#   1: Title "RSTs project TODO list 21/12/2017"
#   2: "Daytime vs. Nighttime RSTs ..."
#   3: "Look at trends on a yearly and monthly basis."
#   4: "Compare classification results between resolutions and models."
#   5: "Create a composite of our classification ..."
#   6: "Remove the polyfit for testing the RST orientation."
#   7: "Send the TODO list." (already struck through)
#
# Apply strikethrough formatting to paragraphs 2 through 6 (both the run
# text and the paragraph mark), leaving paragraph 7 untouched.

for ($i = 2; $i -le 6; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Font.StrikeThrough = 1
}
